$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 57
    3  = 58
    4  = 60
    6  = 62
    7  = 62
    8  = 63
    9  = 65
    10 = 67
    11 = 68
    12 = 69
    13 = 70
    14 = 71
    15 = 72
    16 = 72
    17 = 73
    18 = 76
    20 = 78
    21 = 79
    22 = 80
    23 = 81
    24 = 82
    25 = 82
    26 = 83
    27 = 84
    28 = 85
    29 = 86
    30 = 87
    31 = 87
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
